$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $s = [string]$val
        if ($s -match ", ") {
            $parts = $s -split ", "
            $rev = $parts[($parts.Count - 1)..0]
            $newVal = [string]::Join(", ", $rev)
            $cell.Value2 = $newVal
        }
    }
}
